$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$data = @(
    @(2, 6, 6),
    @(3, 4, 5),
    @(4, 7, 7),
    @(5, 5, 6),
    @(6, 5, 5),
    @(7, 7, 7),
    @(8, 8, 8),
    @(9, 9, 9),
    @(10, 7, 7),
    @(11, 6, 6),
    @(12, 5, 6),
    @(13, 8, 8),
    @(14, 6, 7),
    @(15, 8, 9),
    @(16, 8, 8),
    @(17, 6, 6),
    @(18, 7, 7),
    @(19, 6, 6),
    @(20, 6, 7),
    @(21, 9, 9),
    @(22, 8, 8),
    @(23, 8, 8),
    @(24, 6, 7),
    @(25, 8, 9),
    @(26, 7, 7),
    @(27, 6, 7),
    @(28, 5, 6),
    @(29, 8, 9),
    @(30, 8, 8),
    @(31, 5, 6),
    @(32, 1, 3),
    @(33, 6, 7),
    @(34, 8, 8),
    @(35, 1, 4),
    @(36, 1, 3),
    @(37, 9, 9),
    @(38, 4, 4),
    @(39, 7, 7),
    @(40, 5, 6),
    @(41, 6, 7),
    @(42, 4, 6),
    @(43, 7, 7),
    @(44, 5, 6),
    @(45, 7, 8),
    @(46, 6, 8),
    @(47, 6, 7),
    @(48, 6, 7),
    @(49, 9, 9),
    @(50, 8, 8),
    @(51, 7, 8),
    @(52, 6, 6),
    @(53, 9, 9),
    @(54, 10, 10),
    @(55, 6, 6),
    @(56, 8, 8),
    @(57, 5, 6),
    @(58, 6, 8),
    @(59, 7, 8),
    @(60, 7, 7),
    @(61, 4, 5),
    @(62, 6, 6),
    @(63, 6, 6),
    @(64, 3, 3),
    @(65, 3, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
